$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": add a new data row (14) for 2020-04-XX (serial 43942), pushing
# the two footnote rows down by one.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# Clone the date-cell format from the row above so the new date cell keeps
# the same number format / style as the rest of the column.
$wsAll.Range("A13").Copy()
$wsAll.Range("A14").PasteSpecial(-4122)

$wsAll.Range("A14").Value = 43942
$wsAll.Range("B14").Value = 202
$wsAll.Range("C14").Value = 156
$wsAll.Range("D14").Value = 103
$wsAll.Range("E14").Value = 96
$wsAll.Range("F14").Value = 7
$wsAll.Range("G14").Value = 3
$wsAll.Range("H14").Value = 54

# Shift the two footnote rows down: old row14 -> row15, old row15 -> row16
$wsAll.Range("B15").Value = "※24・34・53・58・59・60・161例目は市外在住者です。"
$wsAll.Range("B14").Copy()
$wsAll.Range("B16").PasteSpecial(-4122)
$wsAll.Range("B16").Value = "※34件調査中"

# ---------------------------------------------------------------------------
# Sheet "kobe": update several running totals, and turn the former footnote
# row (69) into a new data row for 2020-04-XX, adding a fresh footnote row
# (70) below it.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D52").Value = 4
$wsKobe.Range("E52").Value = 43
$wsKobe.Range("D53").Value = 8
$wsKobe.Range("E64").Value = 161
$wsKobe.Range("E65").Value = 173
$wsKobe.Range("E66").Value = 182
$wsKobe.Range("E67").Value = 185
$wsKobe.Range("C68").Value = 1276
$wsKobe.Range("D68").Value = 8
$wsKobe.Range("E68").Value = 193

# Clone the row-68 formatting down into row 69 (new data row). Column B
# (the former footnote cell) already carries the style we want, so it is
# left untouched.
$wsKobe.Range("A68").Copy()
$wsKobe.Range("A69").PasteSpecial(-4122)
$wsKobe.Range("C68:J68").Copy()
$wsKobe.Range("C69:J69").PasteSpecial(-4122)

$wsKobe.Range("A69").Value = 43942
$wsKobe.Range("B69").Value = 0
$wsKobe.Range("C69").Value = 1276
$wsKobe.Range("D69").Value = 9
$wsKobe.Range("E69").Value = 202
$wsKobe.Range("F69").Value = 96
$wsKobe.Range("G69").Value = 90
$wsKobe.Range("H69").Value = 6
$wsKobe.Range("I69").Value = 3
$wsKobe.Range("J69").Value = 50

$wsKobe.Range("B70").Value = "※24・34・53・58・59・60例目は市外在住者です。"

# ---------------------------------------------------------------------------
# Restore the saved cursor position on each sheet (cosmetic, matches the
# "after" sheetViews selection).
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsAll.Activate()
$wsAll.Range("I14").Select()

$wsKobe.Activate()
$wsKobe.Range("C69").Select()

$wsOther.Activate()
$wsOther.Range("H43").Select()

$wsAll.Activate()
